# Update to also include PdeltaAIC as a cov for CG path.
#
# The relation "GR<-det_Clim" moves from row 2 down to row 5 (displaced
# by the newly-surfaced "GR<-Pop_mean" / "GR<-Trait_mean" /
# "Trait_mean<-det_Clim" rows, which now appear first), and its
# coefficient/SD values are recomputed. "Ind_GR<-det_Clim" and
# "Tot_GR<-det_Clim" (rows 6-7) are unaffected.
#
# All of these values look numeric (scientific notation) but must stay
# plain text cells (as in the original file), so we temporarily force a
# text number format while assigning them, then clear the format back
# off again so the cells end up unstyled, exactly like the source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; A = "GR<-Pop_mean";          B = "6.441021e-15"; C = "1.261729e-14"; D = "6.161400e-03" },
    @{ Row = 3; A = "GR<-Trait_mean";        B = "0.000000e+00"; C = "1.047531e-02"; D = "0.000000e+00" },
    @{ Row = 4; A = "Trait_mean<-det_Clim";  B = "4.505846e-02"; C = "2.766086e-02"; D = "4.426837e-01" },
    @{ Row = 5; A = "GR<-det_Clim";          B = "3.476596e-15"; C = "0.000000e+00"; D = "6.167545e-02" },
    @{ Row = 6; A = "Ind_GR<-det_Clim";      B = "0.000000e+00"; C = "0.000000e+00"; D = "0.000000e+00" },
    @{ Row = 7; A = "Tot_GR<-det_Clim";      B = "0.000000e+00"; C = "0.000000e+00"; D = "0.000000e+00" }
)

$rng = $ws.Range("A2:D7")
$rng.NumberFormat = "@"

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.A
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
}

$rng.ClearFormats()
